# Auto-generated script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("site_metrics")
$ws3 = $wb.Worksheets.Item("mk_duration")
$ws4 = $wb.Worksheets.Item("mk_intra_annual")

# ---- site_metrics ----
$ws1.Range("O15").Value = 0.007649406432195554
$ws1.Range("O17").Value = 0.00603356254614023
$ws1.Range("AK17").Value = $True
$ws1.Range("AK20").Value = $True
$ws1.Range("AK21").Value = $True
$ws1.Range("O23").Value = 0.05182091399710077
$ws1.Range("AK23").Value = $True
$ws1.Range("O29").Value = 0.006831240862513246
$ws1.Range("O33").Value = 0.1792874572970606
$ws1.Range("O36").Value = 0.1023922516256202
$ws1.Range("O37").Value = 0.3049926423347166
$ws1.Range("AK42").Value = $True
$ws1.Range("N43").Value = 5.877142857142856
$ws1.Range("O45").Value = 0.1559372094145271
$ws1.Range("N45").Value = 8.478919860627178
$ws1.Range("O57").Value = 0.1713674416333563
$ws1.Range("AK58").Value = $True
$ws1.Range("N59").Value = 9.900595238095237
$ws1.Range("O59").Value = 0.06455196448691045
$ws1.Range("O61").Value = 0.08280397171873011
$ws1.Range("AK62").Value = $True
$ws1.Range("N67").Value = 4.969459141681363
$ws1.Range("O70").Value = 0.0004256360826732732
$ws1.Range("O77").Value = 0.02802196116224578
$ws1.Range("AK78").Value = $True
$ws1.Range("O81").Value = 0.004649956386151363
$ws1.Range("O82").Value = 0.06544405309196726
$ws1.Range("O83").Value = 0.06699524439480588
$ws1.Range("O86").Value = 0.0199608926356963
$ws1.Range("AK88").Value = $True
$ws1.Range("AK91").Value = $True
$ws1.Range("O95").Value = 0.06813302415270088
$ws1.Range("N95").Value = 33.53666666666667
$ws1.Range("Q95").Value = 2.266666666666667
$ws1.Range("AK96").Value = $True
$ws1.Range("AK99").Value = $True
$ws1.Range("O104").Value = 0.04485347931880867
$ws1.Range("O108").Value = 0.4186356387818895
$ws1.Range("O112").Value = 0.2561460881617742
$ws1.Range("N113").Value = 5.739353741496599
$ws1.Range("O113").Value = 0.002990118669602989
$ws1.Range("N115").Value = 3.780320399764845
$ws1.Range("O117").Value = 0.3901976379788569
$ws1.Range("AK119").Value = $True
$ws1.Range("O120").Value = 0.007089722297120988
$ws1.Range("O121").Value = 0.01162218826395471
$ws1.Range("N125").Value = 13.87301587301587
$ws1.Range("AK126").Value = $True
$ws1.Range("O127").Value = 0.0021237068613068
$ws1.Range("AK128").Value = $True
$ws1.Range("AK132").Value = $True
$ws1.Range("AK133").Value = $True
$ws1.Range("AK135").Value = $True
$ws1.Range("AK136").Value = $True
$ws1.Range("N138").Value = 3.270228884590587
$ws1.Range("N142").Value = 7.608333333333334
$ws1.Range("O142").Value = 0.007306643500500288

# ---- mk_duration ----
$ws3.Range("Q4").Value = 12957
$ws3.Range("M4").Value = 0.8605277266625444
$ws3.Range("P4").Value = -21
$ws3.Range("O4").Value = -0.01785714285714286
$ws3.Range("N4").Value = -0.1757024293378821
$ws3.Range("Q18").Value = 814
$ws3.Range("R18").Value = -0.06593406593406594
$ws3.Range("N18").Value = -0.455649764201308
$ws3.Range("S18").Value = 5.393406593406593
$ws3.Range("P18").Value = -14
$ws3.Range("M18").Value = 0.6486418517513601
$ws3.Range("O18").Value = -0.08187134502923976
$ws3.Range("P22").Value = -5
$ws3.Range("K22").Value = "no trend"
$ws3.Range("Q22").Value = 5377.666666666667
$ws3.Range("M22").Value = 0.9565001460423768
$ws3.Range("R22").Value = 0
$ws3.Range("N22").Value = -0.05454601804028707
$ws3.Range("S22").Value = 3.666666666666667
$ws3.Range("O22").Value = -0.007936507936507936
$ws3.Range("L22").Value = $False
$ws3.Range("M26").Value = 0.3577948655088212
$ws3.Range("O26").Value = 0.0858843537414966
$ws3.Range("S26").Value = 4.5
$ws3.Range("N26").Value = 0.9195750583312053
$ws3.Range("P26").Value = 101
$ws3.Range("Q26").Value = 11825.66666666667
$ws3.Range("R29").Value = 0.1818181818181818
$ws3.Range("N29").Value = 2.109233855946642
$ws3.Range("Q29").Value = 2972.666666666667
$ws3.Range("S29").Value = -0.6363636363636362
$ws3.Range("O29").Value = 0.2666666666666667
$ws3.Range("P29").Value = 116
$ws3.Range("M29").Value = 0.03492440092971849
$ws3.Range("O38").Value = -0.07666666666666666
$ws3.Range("P38").Value = -23
$ws3.Range("K38").Value = "no trend"
$ws3.Range("M38").Value = 0.6041785223062051
$ws3.Range("L38").Value = $False
$ws3.Range("S38").Value = 6.566666666666666
$ws3.Range("N38").Value = -0.5184009925889004
$ws3.Range("R38").Value = -0.08888888888888891
$ws3.Range("Q38").Value = 1801
$ws3.Range("O40").Value = -0.09713228492136911
$ws3.Range("M40").Value = 0.3396525330846307
$ws3.Range("N40").Value = -0.9548520263958344
$ws3.Range("R40").Value = -0.0456989247311828
$ws3.Range("Q40").Value = 11863
$ws3.Range("S40").Value = 7.551075268817204
$ws3.Range("P40").Value = -105
$ws3.Range("N43").Value = -0.4309030829119628
$ws3.Range("P43").Value = -42
$ws3.Range("R43").Value = -0.006756756756756757
$ws3.Range("O43").Value = -0.04651162790697674
$ws3.Range("M43").Value = 0.6665388417397997
$ws3.Range("Q43").Value = 9053.333333333334
$ws3.Range("S43").Value = 4.641891891891892
$ws3.Range("O45").Value = 0.0272108843537415
$ws3.Range("P45").Value = 32
$ws3.Range("N45").Value = 0.2679792584830652
$ws3.Range("M45").Value = 0.7887152843391385
$ws3.Range("Q45").Value = 13382
$ws3.Range("L46").Value = $False
$ws3.Range("O46").Value = -0.04873949579831933
$ws3.Range("P46").Value = -29
$ws3.Range("Q46").Value = 4743.666666666667
$ws3.Range("R46").Value = 0
$ws3.Range("K46").Value = "no trend"
$ws3.Range("N46").Value = -0.4065378156740708
$ws3.Range("S46").Value = 11
$ws3.Range("M46").Value = 0.6843474773328535
$ws3.Range("P59").Value = 83
$ws3.Range("N59").Value = 1.274269349956132
$ws3.Range("M59").Value = 0.2025679749526061
$ws3.Range("S59").Value = 5.607142857142858
$ws3.Range("K59").Value = "no trend"
$ws3.Range("Q59").Value = 4141
$ws3.Range("O59").Value = 0.1571969696969697
$ws3.Range("R59").Value = 0.1495535714285714
$ws3.Range("L59").Value = $False
$ws3.Range("S66").Value = 4.384615384615385
$ws3.Range("Q66").Value = 7905.666666666667
$ws3.Range("O66").Value = -0.06707317073170732
$ws3.Range("P66").Value = -55
$ws3.Range("N66").Value = -0.6073296868822874
$ws3.Range("M66").Value = 0.5436321407219808
$ws3.Range("R66").Value = -0.01923076923076923
$ws3.Range("P67").Value = 1
$ws3.Range("K67").Value = "no trend"
$ws3.Range("O67").Value = 0.002463054187192118
$ws3.Range("Q67").Value = 2837
$ws3.Range("R67").Value = 0
$ws3.Range("M67").Value = 1
$ws3.Range("N67").Value = 0
$ws3.Range("S67").Value = 4.75
$ws3.Range("L67").Value = $False
$ws3.Range("K70").Value = "no trend"
$ws3.Range("P70").Value = 33
$ws3.Range("R70").Value = 0.0625
$ws3.Range("M70").Value = 0.3949411107763483
$ws3.Range("Q70").Value = 1415
$ws3.Range("O70").Value = 0.1304347826086956
$ws3.Range("S70").Value = 2.3125
$ws3.Range("N70").Value = 0.8506908437621611
$ws3.Range("L70").Value = $False
$ws3.Range("L72").Value = $False
$ws3.Range("K72").Value = "no trend"
$ws3.Range("S72").Value = 1.666666666666667
$ws3.Range("O72").Value = 0.1660079051383399
$ws3.Range("Q72").Value = 1308.666666666667
$ws3.Range("N72").Value = 1.133363798497767
$ws3.Range("P72").Value = 42
$ws3.Range("M72").Value = 0.2570615100989806
$ws3.Range("F95").Value = 0.06666666666666667
$ws3.Range("D95").Value = 0.7665252995235718
$ws3.Range("E95").Value = 0.2969229955832361
$ws3.Range("G95").Value = 7
$ws3.Range("P95").Value = -41
$ws3.Range("O95").Value = -0.0553306342780027
$ws3.Range("N95").Value = -0.5542562584220408
$ws3.Range("M95").Value = 0.5794034813424265
$ws3.Range("P103").Value = -7
$ws3.Range("M103").Value = 0.9308844341393401
$ws3.Range("N103").Value = -0.08673213308848518
$ws3.Range("O103").Value = -0.01176470588235294
$ws3.Range("Q103").Value = 4785.666666666667
$ws3.Range("S103").Value = 3.4
$ws3.Range("Q106").Value = 1884.666666666667
$ws3.Range("M106").Value = 0.9816225867151522
$ws3.Range("O106").Value = -0.006153846153846154
$ws3.Range("N106").Value = -0.02303470874272992
$ws3.Range("S106").Value = 1
$ws3.Range("P106").Value = -2
$ws3.Range("M113").Value = 0.4539693793700228
$ws3.Range("S113").Value = 2.523199023199023
$ws3.Range("P113").Value = 77
$ws3.Range("O113").Value = 0.07777777777777778
$ws3.Range("N113").Value = 0.7488139020976839
$ws3.Range("R113").Value = 0.02167277167277167
$ws3.Range("Q113").Value = 10301
$ws3.Range("Q115").Value = 3448.333333333333
$ws3.Range("N115").Value = 0.8174037338950666
$ws3.Range("P115").Value = 49
$ws3.Range("O115").Value = 0.1053763440860215
$ws3.Range("R115").Value = 0.03571428571428581
$ws3.Range("M115").Value = 0.4136977425759303
$ws3.Range("S115").Value = 2.019841269841268
$ws3.Range("R118").Value = -0.07142857142857141
$ws3.Range("S118").Value = 6.571428571428571
$ws3.Range("P118").Value = -3
$ws3.Range("Q118").Value = 588.3333333333334
$ws3.Range("L118").Value = $False
$ws3.Range("K118").Value = "no trend"
$ws3.Range("M118").Value = 0.9342847091688518
$ws3.Range("O118").Value = -0.02205882352941177
$ws3.Range("N118").Value = -0.08245524152850489
$ws3.Range("S120").Value = 7.6125
$ws3.Range("O120").Value = -0.1057471264367816
$ws3.Range("R120").Value = -0.125
$ws3.Range("Q120").Value = 3108.666666666667
$ws3.Range("L120").Value = $False
$ws3.Range("P120").Value = -46
$ws3.Range("N120").Value = -0.8070964472117726
$ws3.Range("M120").Value = 0.4196109188429
$ws3.Range("K120").Value = "no trend"
$ws3.Range("Q122").Value = 8904
$ws3.Range("O122").Value = -0.0664451827242525
$ws3.Range("P122").Value = -60
$ws3.Range("N122").Value = -0.6252582575139637
$ws3.Range("S122").Value = 5
$ws3.Range("M122").Value = 0.5318015717062972
$ws3.Range("O125").Value = -0.2044334975369458
$ws3.Range("R125").Value = -0.3674749163879598
$ws3.Range("Q125").Value = 2773.666666666667
$ws3.Range("N125").Value = -1.556993210399037
$ws3.Range("M125").Value = 0.1194720974217798
$ws3.Range("P125").Value = -83
$ws3.Range("S125").Value = 10.89464882943144
$ws3.Range("M126").Value = 0.9241325721763831
$ws3.Range("P126").Value = -11
$ws3.Range("O126").Value = -0.0106280193236715
$ws3.Range("Q126").Value = 11027
$ws3.Range("S126").Value = 7
$ws3.Range("N126").Value = -0.09522945803722153
$ws3.Range("S127").Value = 4.752380952380953
$ws3.Range("Q127").Value = 1789
$ws3.Range("O127").Value = 0.1966666666666667
$ws3.Range("R127").Value = 0.2984126984126984
$ws3.Range("M127").Value = 0.1702909531150576
$ws3.Range("N127").Value = 1.371269521507303
$ws3.Range("N138").Value = -0.1010465043413376
$ws3.Range("S138").Value = 3
$ws3.Range("R138").Value = 0
$ws3.Range("Q138").Value = 11850.66666666667
$ws3.Range("P138").Value = -12
$ws3.Range("M138").Value = 0.9195135439406887
$ws3.Range("O138").Value = -0.01110083256244218
$ws3.Range("N142").Value = -0.7625640231279325
$ws3.Range("Q142").Value = 5785
$ws3.Range("R142").Value = -0.05555555555555561
$ws3.Range("S142").Value = 6.200000000000001
$ws3.Range("O142").Value = -0.08858858858858859
$ws3.Range("P142").Value = -59
$ws3.Range("M142").Value = 0.4457234473481497

# ---- mk_intra_annual ----
$ws4.Range("Q4").Value = 12194.66666666667
$ws4.Range("M4").Value = 0.7104306898553399
$ws4.Range("P4").Value = 42
$ws4.Range("O4").Value = 0.03571428571428571
$ws4.Range("N4").Value = 0.3712777212457833
$ws4.Range("Q18").Value = 793
$ws4.Range("N18").Value = -0.142044164845687
$ws4.Range("S18").Value = 3
$ws4.Range("P18").Value = -5
$ws4.Range("M18").Value = 0.8870451208364276
$ws4.Range("O18").Value = -0.02923976608187134
$ws4.Range("P22").Value = -66
$ws4.Range("K22").Value = "no trend"
$ws4.Range("Q22").Value = 5284.666666666667
$ws4.Range("M22").Value = 0.3712480546828387
$ws4.Range("R22").Value = -0.0303030303030303
$ws4.Range("N22").Value = -0.8941380096891555
$ws4.Range("S22").Value = 4.53030303030303
$ws4.Range("O22").Value = -0.1047619047619048
$ws4.Range("L22").Value = $False
$ws4.Range("M26").Value = 0.9850601489981838
$ws4.Range("O26").Value = -0.002551020408163265
$ws4.Range("S26").Value = 1
$ws4.Range("N26").Value = -0.01872542073030786
$ws4.Range("P26").Value = -3
$ws4.Range("Q26").Value = 11407.66666666667
$ws4.Range("N29").Value = 1.125392299802487
$ws4.Range("Q29").Value = 2938
$ws4.Range("S29").Value = 1.5
$ws4.Range("O29").Value = 0.1425287356321839
$ws4.Range("P29").Value = 62
$ws4.Range("M29").Value = 0.2604228326234921
$ws4.Range("O38").Value = 0.05666666666666666
$ws4.Range("P38").Value = 17
$ws4.Range("K38").Value = "no trend"
$ws4.Range("M38").Value = 0.700283863420285
$ws4.Range("L38").Value = $False
$ws4.Range("S38").Value = 2
$ws4.Range("N38").Value = 0.3849373086914444
$ws4.Range("Q38").Value = 1727.666666666667
$ws4.Range("O40").Value = -0.06197964847363552
$ws4.Range("M40").Value = 0.5294927145023132
$ws4.Range("N40").Value = -0.6287805816154457
$ws4.Range("Q40").Value = 11017.66666666667
$ws4.Range("P40").Value = -67
$ws4.Range("N43").Value = -0.3950653277746584
$ws4.Range("P43").Value = -38
$ws4.Range("O43").Value = -0.04208194905869325
$ws4.Range("M43").Value = 0.6927946763304649
$ws4.Range("Q43").Value = 8771.333333333334
$ws4.Range("O45").Value = -0.102891156462585
$ws4.Range("P45").Value = -121
$ws4.Range("N45").Value = -1.05830640392219
$ws4.Range("M45").Value = 0.2899157743181568
$ws4.Range("Q45").Value = 12857
$ws4.Range("L46").Value = $False
$ws4.Range("O46").Value = -0.146218487394958
$ws4.Range("P46").Value = -87
$ws4.Range("Q46").Value = 4139.666666666667
$ws4.Range("K46").Value = "no trend"
$ws4.Range("N46").Value = -1.33664403637402
$ws4.Range("S46").Value = 1
$ws4.Range("M46").Value = 0.1813388574991128
$ws4.Range("P59").Value = 86
$ws4.Range("N59").Value = 1.388785070569616
$ws4.Range("M59").Value = 0.1648981163691001
$ws4.Range("S59").Value = 2
$ws4.Range("K59").Value = "no trend"
$ws4.Range("Q59").Value = 3746
$ws4.Range("O59").Value = 0.1628787878787879
$ws4.Range("R59").Value = 0
$ws4.Range("L59").Value = $False
$ws4.Range("S66").Value = 3
$ws4.Range("Q66").Value = 7586.333333333333
$ws4.Range("O66").Value = 0.006097560975609756
$ws4.Range("P66").Value = 5
$ws4.Range("N66").Value = 0.04592445708074954
$ws4.Range("M66").Value = 0.9633704608077815
$ws4.Range("R66").Value = 0
$ws4.Range("P67").Value = -1
$ws4.Range("K67").Value = "no trend"
$ws4.Range("O67").Value = -0.002463054187192118
$ws4.Range("Q67").Value = 2757
$ws4.Range("R67").Value = 0
$ws4.Range("M67").Value = 1
$ws4.Range("N67").Value = 0
$ws4.Range("S67").Value = 3
$ws4.Range("L67").Value = $False
$ws4.Range("K70").Value = "no trend"
$ws4.Range("P70").Value = -33
$ws4.Range("M70").Value = 0.3869973371978499
$ws4.Range("Q70").Value = 1368.333333333333
$ws4.Range("O70").Value = -0.1304347826086956
$ws4.Range("S70").Value = 2
$ws4.Range("N70").Value = -0.8650755249251965
$ws4.Range("L70").Value = $False
$ws4.Range("L72").Value = $False
$ws4.Range("K72").Value = "no trend"
$ws4.Range("S72").Value = 1
$ws4.Range("O72").Value = -0.04347826086956522
$ws4.Range("Q72").Value = 1287.666666666667
$ws4.Range("N72").Value = -0.2786751690725666
$ws4.Range("P72").Value = -11
$ws4.Range("M72").Value = 0.7804941200970492
$ws4.Range("F95").Value = -0.2380952380952381
$ws4.Range("H95").Value = 314.3333333333333
$ws4.Range("D95").Value = 0.1758384813657532
$ws4.Range("E95").Value = -1.353680031929301
$ws4.Range("G95").Value = -25
$ws4.Range("P95").Value = -73
$ws4.Range("O95").Value = -0.09851551956815115
$ws4.Range("N95").Value = -1.006787885851035
$ws4.Range("Q95").Value = 5114.333333333333
$ws4.Range("M95").Value = 0.3140367173593073
$ws4.Range("P103").Value = -93
$ws4.Range("M103").Value = 0.1797484873874016
$ws4.Range("N103").Value = -1.341529842258298
$ws4.Range("O103").Value = -0.1563025210084034
$ws4.Range("Q103").Value = 4703
$ws4.Range("K106").Value = "no trend"
$ws4.Range("Q106").Value = 1846.333333333333
$ws4.Range("M106").Value = 0.485065843171268
$ws4.Range("N106").Value = 0.6981780638018836
$ws4.Range("S106").Value = 1
$ws4.Range("L106").Value = $False
$ws4.Range("P106").Value = 31
$ws4.Range("O106").Value = 0.09538461538461539
$ws4.Range("M113").Value = 0.3127371169525379
$ws4.Range("P113").Value = -102
$ws4.Range("O113").Value = -0.103030303030303
$ws4.Range("N113").Value = -1.009495378434651
$ws4.Range("Q113").Value = 10010
$ws4.Range("Q115").Value = 3386
$ws4.Range("N115").Value = -1.323266287655244
$ws4.Range("P115").Value = -78
$ws4.Range("O115").Value = -0.167741935483871
$ws4.Range("R115").Value = -0.08333333333333333
$ws4.Range("M115").Value = 0.1857468386547563
$ws4.Range("S115").Value = 5.25
$ws4.Range("S118").Value = 2
$ws4.Range("P118").Value = 1
$ws4.Range("Q118").Value = 534.3333333333334
$ws4.Range("L118").Value = $False
$ws4.Range("K118").Value = "no trend"
$ws4.Range("M118").Value = 1
$ws4.Range("O118").Value = 0.007352941176470588
$ws4.Range("N118").Value = 0
$ws4.Range("S120").Value = 2
$ws4.Range("O120").Value = -0.06436781609195402
$ws4.Range("Q120").Value = 2972.666666666667
$ws4.Range("L120").Value = $False
$ws4.Range("P120").Value = -28
$ws4.Range("N120").Value = -0.4952114270483421
$ws4.Range("M120").Value = 0.6204508852729831
$ws4.Range("K120").Value = "no trend"
$ws4.Range("Q122").Value = 8527.666666666666
$ws4.Range("O122").Value = -0.001107419712070875
$ws4.Range("P122").Value = -1
$ws4.Range("N122").Value = 0
$ws4.Range("M122").Value = 1
$ws4.Range("O125").Value = 0.07389162561576355
$ws4.Range("Q125").Value = 2646.666666666667
$ws4.Range("N125").Value = 0.5637004403031675
$ws4.Range("M125").Value = 0.5729580176967639
$ws4.Range("P125").Value = 30
$ws4.Range("S125").Value = 2
$ws4.Range("M126").Value = 0.755851007513781
$ws4.Range("P126").Value = 33
$ws4.Range("O126").Value = 0.03188405797101449
$ws4.Range("Q126").Value = 10591.66666666667
$ws4.Range("N126").Value = 0.3109337223606067
$ws4.Range("S127").Value = 1
$ws4.Range("Q127").Value = 1715.666666666667
$ws4.Range("O127").Value = 0.15
$ws4.Range("M127").Value = 0.2881116829209884
$ws4.Range("N127").Value = 1.062273187484946
$ws4.Range("N138").Value = 1.270772113468216
$ws4.Range("S138").Value = 4.378378378378378
$ws4.Range("R138").Value = 0.02702702702702703
$ws4.Range("Q138").Value = 11622.66666666667
$ws4.Range("P138").Value = 138
$ws4.Range("K138").Value = "no trend"
$ws4.Range("L138").Value = $False
$ws4.Range("M138").Value = 0.2038097322334718
$ws4.Range("O138").Value = 0.1276595744680851
$ws4.Range("N142").Value = 0.9204637578316442
$ws4.Range("Q142").Value = 5619.333333333333
$ws4.Range("R142").Value = 0
$ws4.Range("S142").Value = 3
$ws4.Range("L142").Value = $False
$ws4.Range("K142").Value = "no trend"
$ws4.Range("O142").Value = 0.1051051051051051
$ws4.Range("P142").Value = 70
$ws4.Range("M142").Value = 0.35733046379264
